$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing Tuesday (column E) hours for the week of row 8
$ws.Range("E8").Value = 6.5

# Update the active selection to match the author's final cursor position
$ws.Range("G15").Select()
